$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: May - update 2021 (H) value
$ws.Range("H6").Value = 108

# Row 9: August - update label to reflect new "through" date, and update values
$ws.Range("A9").Value = "August (through 08-25)"
$ws.Range("C9").Value = 59
$ws.Range("D9").Value = 71
$ws.Range("E9").Value = 46
$ws.Range("F9").Value = 37
$ws.Range("G9").Value = 141
$ws.Range("H9").Value = 124

# Row 10: Total - update values
$ws.Range("C10").Value = 361
$ws.Range("D10").Value = 536
$ws.Range("E10").Value = 471
$ws.Range("F10").Value = 341
$ws.Range("G10").Value = 762
$ws.Range("H10").Value = 1038
